$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("100_1")

$ws.Range("B11").Value = 1
$ws.Range("B21").Value = 1
$ws.Range("B26").Value = 1
